$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 63: shift betting-market data (teams/odds/timestamps/url) per scraper re-run
$ws.Cells.Item(63,6).Value = 'Richards Bay'
$ws.Cells.Item(63,7).Value = 0
$ws.Cells.Item(63,8).Value = 'Polokwane'
$ws.Cells.Item(63,9).Value = 1
$ws.Cells.Item(63,10).Value = 2.58
$ws.Cells.Item(63,11).Value = '03/10/2023 01:12'
$ws.Cells.Item(63,12).Value = 2.42
$ws.Cells.Item(63,13).Value = '04/10/2023 19:21'
$ws.Cells.Item(63,14).Value = 2.85
$ws.Cells.Item(63,15).Value = '03/10/2023 01:12'
$ws.Cells.Item(63,16).Value = 2.83
$ws.Cells.Item(63,17).Value = '04/10/2023 19:21'
$ws.Cells.Item(63,18).Value = 3.17
$ws.Cells.Item(63,19).Value = '03/10/2023 01:12'
$ws.Cells.Item(63,20).Value = 3.61
$ws.Cells.Item(63,21).Value = '04/10/2023 19:21'
$ws.Cells.Item(63,22).Value = 'https://www.betexplorer.com/football/south-africa/premier-league/richards-bay-polokwane-city/4SkrPj8b/'

# Row 64: shift betting-market data (teams/odds/timestamps/url) per scraper re-run
$ws.Cells.Item(64,6).Value = 'Stellenbosch'
$ws.Cells.Item(64,7).Value = 3
$ws.Cells.Item(64,8).Value = 'TS Galaxy'
$ws.Cells.Item(64,9).Value = 0
$ws.Cells.Item(64,10).Value = 2.27
$ws.Cells.Item(64,11).Value = '03/10/2023 01:12'
$ws.Cells.Item(64,12).Value = 2.3
$ws.Cells.Item(64,13).Value = '04/10/2023 19:22'
$ws.Cells.Item(64,14).Value = 2.99
$ws.Cells.Item(64,15).Value = '03/10/2023 01:12'
$ws.Cells.Item(64,16).Value = 2.93
$ws.Cells.Item(64,17).Value = '04/10/2023 19:22'
$ws.Cells.Item(64,18).Value = 3.58
$ws.Cells.Item(64,19).Value = '03/10/2023 01:12'
$ws.Cells.Item(64,20).Value = 3.74
$ws.Cells.Item(64,21).Value = '04/10/2023 19:22'
$ws.Cells.Item(64,22).Value = 'https://www.betexplorer.com/football/south-africa/premier-league/stellenbosch-fc-ts-galaxy/niumOAN4/'

# Row 65: shift betting-market data (teams/odds/timestamps/url) per scraper re-run
$ws.Cells.Item(65,6).Value = 'Cape Town Spurs'
$ws.Cells.Item(65,7).Value = 1
$ws.Cells.Item(65,8).Value = 'Supersport Utd'
$ws.Cells.Item(65,9).Value = 2
$ws.Cells.Item(65,10).Value = 4.52
$ws.Cells.Item(65,11).Value = '30/09/2023 14:13'
$ws.Cells.Item(65,12).Value = 5.17
$ws.Cells.Item(65,13).Value = '04/10/2023 19:21'
$ws.Cells.Item(65,14).Value = 3.07
$ws.Cells.Item(65,15).Value = '30/09/2023 14:13'
$ws.Cells.Item(65,16).Value = 3.52
$ws.Cells.Item(65,17).Value = '04/10/2023 19:21'
$ws.Cells.Item(65,18).Value = 1.96
$ws.Cells.Item(65,19).Value = '30/09/2023 14:13'
$ws.Cells.Item(65,20).Value = 1.75
$ws.Cells.Item(65,21).Value = '04/10/2023 19:21'
$ws.Cells.Item(65,22).Value = 'https://www.betexplorer.com/football/south-africa/premier-league/cape-town-spurs-supersport-utd/fHreMlhH/'

# Row 75: shift betting-market data (teams/odds/timestamps/url) per scraper re-run
$ws.Cells.Item(75,6).Value = 'Stellenbosch'
$ws.Cells.Item(75,7).Value = 2
$ws.Cells.Item(75,8).Value = 'Richards Bay'
$ws.Cells.Item(75,9).Value = 1
$ws.Cells.Item(75,10).Value = 2.03
$ws.Cells.Item(75,11).Value = '06/11/2023 15:19'
$ws.Cells.Item(75,12).Value = 1.99
$ws.Cells.Item(75,13).Value = '07/11/2023 18:22'
$ws.Cells.Item(75,14).Value = 3.24
$ws.Cells.Item(75,15).Value = '06/11/2023 15:19'
$ws.Cells.Item(75,16).Value = 3.17
$ws.Cells.Item(75,17).Value = '07/11/2023 18:22'
$ws.Cells.Item(75,18).Value = 4.1
$ws.Cells.Item(75,19).Value = '06/11/2023 15:19'
$ws.Cells.Item(75,20).Value = 4.37
$ws.Cells.Item(75,21).Value = '07/11/2023 18:21'
$ws.Cells.Item(75,22).Value = 'https://www.betexplorer.com/football/south-africa/premier-league/stellenbosch-fc-richards-bay/hGIEFwy2/'

# Row 76: shift betting-market data (teams/odds/timestamps/url) per scraper re-run
$ws.Cells.Item(76,6).Value = 'Royal AM'
$ws.Cells.Item(76,7).Value = 1
$ws.Cells.Item(76,8).Value = 'Golden Arrows'
$ws.Cells.Item(76,9).Value = 0
$ws.Cells.Item(76,10).Value = 3.64
$ws.Cells.Item(76,11).Value = '06/11/2023 15:19'
$ws.Cells.Item(76,12).Value = 3.25
$ws.Cells.Item(76,13).Value = '07/11/2023 18:26'
$ws.Cells.Item(76,14).Value = 3.19
$ws.Cells.Item(76,15).Value = '06/11/2023 15:19'
$ws.Cells.Item(76,16).Value = 3.07
$ws.Cells.Item(76,17).Value = '07/11/2023 18:26'
$ws.Cells.Item(76,18).Value = 2.19
$ws.Cells.Item(76,19).Value = '06/11/2023 15:19'
$ws.Cells.Item(76,20).Value = 2.43
$ws.Cells.Item(76,21).Value = '07/11/2023 18:26'
$ws.Cells.Item(76,22).Value = 'https://www.betexplorer.com/football/south-africa/premier-league/royal-am-golden-arrows/EPJAGJLe/'

# Row 77: shift betting-market data (teams/odds/timestamps/url) per scraper re-run
$ws.Cells.Item(77,6).Value = 'Cape Town City'
$ws.Cells.Item(77,7).Value = 2
$ws.Cells.Item(77,8).Value = 'Chippa Utd.'
$ws.Cells.Item(77,9).Value = 1
$ws.Cells.Item(77,10).Value = 1.71
$ws.Cells.Item(77,11).Value = '06/11/2023 15:19'
$ws.Cells.Item(77,12).Value = 1.83
$ws.Cells.Item(77,13).Value = '07/11/2023 18:28'
$ws.Cells.Item(77,14).Value = 3.53
$ws.Cells.Item(77,15).Value = '06/11/2023 15:19'
$ws.Cells.Item(77,16).Value = 3.19
$ws.Cells.Item(77,17).Value = '07/11/2023 18:28'
$ws.Cells.Item(77,18).Value = 5.54
$ws.Cells.Item(77,19).Value = '06/11/2023 15:19'
$ws.Cells.Item(77,20).Value = 5.37
$ws.Cells.Item(77,21).Value = '07/11/2023 18:28'
$ws.Cells.Item(77,22).Value = 'https://www.betexplorer.com/football/south-africa/premier-league/cape-town-city-chippa-utd/S8L2Iuiq/'

# Row 78: shift betting-market data (teams/odds/timestamps/url) per scraper re-run
$ws.Cells.Item(78,6).Value = 'Orlando Pirates'
$ws.Cells.Item(78,7).Value = 1
$ws.Cells.Item(78,8).Value = 'Sekhukhune'
$ws.Cells.Item(78,9).Value = 0
$ws.Cells.Item(78,10).Value = 1.62
$ws.Cells.Item(78,11).Value = '06/11/2023 13:00'
$ws.Cells.Item(78,12).Value = 1.66
$ws.Cells.Item(78,13).Value = '07/11/2023 18:22'
$ws.Cells.Item(78,14).Value = 3.8
$ws.Cells.Item(78,15).Value = '06/11/2023 13:00'
$ws.Cells.Item(78,16).Value = 3.55
$ws.Cells.Item(78,17).Value = '07/11/2023 18:22'
$ws.Cells.Item(78,18).Value = 5.04
$ws.Cells.Item(78,19).Value = '06/11/2023 13:00'
$ws.Cells.Item(78,20).Value = 6.06
$ws.Cells.Item(78,21).Value = '07/11/2023 18:22'
$ws.Cells.Item(78,22).Value = 'https://www.betexplorer.com/football/south-africa/premier-league/orlando-pirates-sekhukhune/6aK6Ha6k/'

# Row 83: shift betting-market data (teams/odds/timestamps/url) per scraper re-run
$ws.Cells.Item(83,6).Value = 'Kaizer Chiefs'
$ws.Cells.Item(83,7).Value = 0
$ws.Cells.Item(83,8).Value = 'Orlando Pirates'
$ws.Cells.Item(83,9).Value = 1
$ws.Cells.Item(83,10).Value = 2.79
$ws.Cells.Item(83,11).Value = '11/11/2023 07:48'
$ws.Cells.Item(83,12).Value = 3.17
$ws.Cells.Item(83,13).Value = '11/11/2023 14:21'
$ws.Cells.Item(83,14).Value = 2.84
$ws.Cells.Item(83,15).Value = '11/11/2023 07:48'
$ws.Cells.Item(83,16).Value = 2.87
$ws.Cells.Item(83,17).Value = '11/11/2023 14:21'
$ws.Cells.Item(83,18).Value = 2.79
$ws.Cells.Item(83,19).Value = '11/11/2023 07:48'
$ws.Cells.Item(83,20).Value = 2.64
$ws.Cells.Item(83,21).Value = '11/11/2023 14:21'
$ws.Cells.Item(83,22).Value = 'https://www.betexplorer.com/football/south-africa/premier-league/kaizer-chiefs-orlando-pirates/0MutWbLr/'

# Row 84: shift betting-market data (teams/odds/timestamps/url) per scraper re-run
$ws.Cells.Item(84,6).Value = 'Polokwane'
$ws.Cells.Item(84,7).Value = 0
$ws.Cells.Item(84,8).Value = 'Swallows'
$ws.Cells.Item(84,9).Value = 0
$ws.Cells.Item(84,10).Value = 2.65
$ws.Cells.Item(84,11).Value = '11/11/2023 07:48'
$ws.Cells.Item(84,12).Value = 2.76
$ws.Cells.Item(84,13).Value = '11/11/2023 14:21'
$ws.Cells.Item(84,14).Value = 2.86
$ws.Cells.Item(84,15).Value = '11/11/2023 07:48'
$ws.Cells.Item(84,16).Value = 2.84
$ws.Cells.Item(84,17).Value = '11/11/2023 14:21'
$ws.Cells.Item(84,18).Value = 2.92
$ws.Cells.Item(84,19).Value = '11/11/2023 07:48'
$ws.Cells.Item(84,20).Value = 3.04
$ws.Cells.Item(84,21).Value = '11/11/2023 14:21'
$ws.Cells.Item(84,22).Value = 'https://www.betexplorer.com/football/south-africa/premier-league/polokwane-city-swallows-fc/fqupVIzk/'

# Row 92: shift betting-market data (teams/odds/timestamps/url) per scraper re-run
$ws.Cells.Item(92,6).Value = 'TS Galaxy'
$ws.Cells.Item(92,7).Value = 3
$ws.Cells.Item(92,8).Value = 'Polokwane'
$ws.Cells.Item(92,9).Value = 0
$ws.Cells.Item(92,10).Value = 2.84
$ws.Cells.Item(92,11).Value = '18/11/2023 14:42'
$ws.Cells.Item(92,12).Value = 2.65
$ws.Cells.Item(92,13).Value = '26/11/2023 14:26'
$ws.Cells.Item(92,14).Value = 2.75
$ws.Cells.Item(92,15).Value = '18/11/2023 14:42'
$ws.Cells.Item(92,16).Value = 2.82
$ws.Cells.Item(92,17).Value = '26/11/2023 14:26'
$ws.Cells.Item(92,18).Value = 2.96
$ws.Cells.Item(92,19).Value = '18/11/2023 14:42'
$ws.Cells.Item(92,20).Value = 3.21
$ws.Cells.Item(92,21).Value = '26/11/2023 14:26'
$ws.Cells.Item(92,22).Value = 'https://www.betexplorer.com/football/south-africa/premier-league/ts-galaxy-polokwane-city/nDYRuzk7/'

# Row 93: shift betting-market data (teams/odds/timestamps/url) per scraper re-run
$ws.Cells.Item(93,6).Value = 'Swallows'
$ws.Cells.Item(93,7).Value = 0
$ws.Cells.Item(93,8).Value = 'Kaizer Chiefs'
$ws.Cells.Item(93,9).Value = 1
$ws.Cells.Item(93,10).Value = 2.54
$ws.Cells.Item(93,11).Value = '18/11/2023 14:42'
$ws.Cells.Item(93,12).Value = 2.76
$ws.Cells.Item(93,13).Value = '26/11/2023 14:26'
$ws.Cells.Item(93,14).Value = 2.91
$ws.Cells.Item(93,15).Value = '18/11/2023 14:42'
$ws.Cells.Item(93,16).Value = 2.97
$ws.Cells.Item(93,17).Value = '26/11/2023 14:26'
$ws.Cells.Item(93,18).Value = 3.16
$ws.Cells.Item(93,19).Value = '18/11/2023 14:42'
$ws.Cells.Item(93,20).Value = 2.9
$ws.Cells.Item(93,21).Value = '26/11/2023 14:26'
$ws.Cells.Item(93,22).Value = 'https://www.betexplorer.com/football/south-africa/premier-league/swallows-fc-kaizer-chiefs/hjZNtGZ0/'

# New row 95: newly scraped fixture (Cape Town Spurs vs Stellenbosch)
$ws.Range("A94:V94").Copy()
$ws.Range("A95:V95").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(95,1).Value = 94
$ws.Cells.Item(95,2).Value = 'south-africa'
$ws.Cells.Item(95,3).Value = 'premier-league'
$ws.Cells.Item(95,4).Value = '2023-2024'
$ws.Cells.Item(95,5).Value = 45258.77083333334
$ws.Cells.Item(95,6).Value = 'Cape Town Spurs'
$ws.Cells.Item(95,7).Value = 0
$ws.Cells.Item(95,8).Value = 'Stellenbosch'
$ws.Cells.Item(95,9).Value = 3
$ws.Cells.Item(95,10).Value = 3.2
$ws.Cells.Item(95,11).Value = '25/11/2023 18:13'
$ws.Cells.Item(95,12).Value = 4.02
$ws.Cells.Item(95,13).Value = '28/11/2023 18:23'
$ws.Cells.Item(95,14).Value = 2.92
$ws.Cells.Item(95,15).Value = '25/11/2023 18:13'
$ws.Cells.Item(95,16).Value = 3.13
$ws.Cells.Item(95,17).Value = '28/11/2023 18:23'
$ws.Cells.Item(95,18).Value = 2.51
$ws.Cells.Item(95,19).Value = '25/11/2023 18:13'
$ws.Cells.Item(95,20).Value = 2.09
$ws.Cells.Item(95,21).Value = '28/11/2023 18:23'
$ws.Cells.Item(95,22).Value = 'https://www.betexplorer.com/football/south-africa/premier-league/cape-town-spurs-stellenbosch-fc/QuaX20id/'
